# Add the new "ManageClass" sheet (ManageClass test data for the LMS suite)
# after the existing "AddNewUserDetails" sheet, and populate it with the
# batch / class management header + sample row.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ManageClass"

# ---- Header row -----------------------------------------------------------
$ws.Range("A1").Value = "BatchName"
$ws.Range("B1").Value = "Description"
$ws.Range("C1").Value = "ProgramName"
$ws.Range("D1").Value = "NoOfClasses"
$ws.Range("E1").Value = "ValidBatchId"
$ws.Range("F1").Value = "InvalidBatchId"
$ws.Range("G1").Value = "ValidClassNo"
$ws.Range("H1").Value = "InvalidClassNo"
$ws.Range("I1").Value = "ValidClassDate"
$ws.Range("J1").Value = "InvalidClassDate"
$ws.Range("K1").Value = "ValidClassTopic"
$ws.Range("L1").Value = "InvalidClassTopic"
$ws.Range("M1").Value = "ValidStaffId"
$ws.Range("N1").Value = "InvalidStaffId"
$ws.Range("O1").Value = "ClassDescription"
$ws.Range("P1").Value = "Comments"
$ws.Range("Q1").Value = "Notes"
$ws.Range("R1").Value = "Recordings"

# ---- Data row ---------------------------------------------------------------
$ws.Range("A2").Value = "SudhaBatch"
$ws.Range("B2").Value = "SDET+JAVA"
$ws.Range("C2").Value = "QA-Automation"

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '"14"'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '"1343"'

$ws.Range("F2").Value = '"327654"'
$ws.Range("G2").Value = '"23"'
$ws.Range("H2").Value = '"6897"'

$ws.Range("I2").NumberFormat = "mm-dd-yy"
$ws.Range("I2").Value = '"12/22/2023"'

$ws.Range("J2").NumberFormat = "d-mmm-yy"
$ws.Range("J2").Value = '"9/22/2024"'
$ws.Range("K2").NumberFormat = "d-mmm-yy"
$ws.Range("K2").Value = "Python"
$ws.Range("L2").NumberFormat = "d-mmm-yy"
$ws.Range("L2").Value = "Science"

$ws.Range("M2").Value = '"U78"'
$ws.Range("N2").Value = '"U90000"'
$ws.Range("O2").Value = "Advanced Python"
$ws.Range("P2").Value = "SelfLearning"
$ws.Range("Q2").Value = "python Notes"
$ws.Range("R2").Value = 'c:\\Recordings'

$ws.Range("A1:R2").Select() | Out-Null
